$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Insert()
$ws.Range("D7").Value = 43465
$v = $ws.Range("D7").Value()
Write-Host "D7:" $v
$v2 = $ws.Range("E7").Value()
Write-Host "E7 (should be old D 43100):" $v2

# Check style of D7 vs E7 (should both be same style as old D, i.e. style index 2/date format)
Write-Host "D7 NumberFormat:" $ws.Range("D7").NumberFormat
Write-Host "E7 NumberFormat:" $ws.Range("E7").NumberFormat
